$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in row 2
$ws.Range("E2").Value = 6
$ws.Range("G2").Value = -3
$ws.Range("H2").Value = 13

# Update the active selection from I1 to E2
$ws.Range("E2").Select()
